# Add a new "24-nov" column (CN) to the right of the existing "23-nov"
# column (CM), mirroring the layout/format of column CM, and fill in the
# per-product counts for that new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's formatting (text format, style index 1) from CM1
# onto the new header cell CN1, then set its value.
$ws.Range("CM1").Copy()
$ws.Range("CN1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("CN1").Value = "24-nov"

# Copy the data cells' formatting (centered integer format, style index 2)
# from column CM onto column CN for every data row, then set the new
# counts for 24-nov.
$ws.Range("CM2:CM11").Copy()
$ws.Range("CN2:CN11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("CN2").Value = 13
$ws.Range("CN3").Value = 9
$ws.Range("CN4").Value = 8
$ws.Range("CN5").Value = 12
$ws.Range("CN6").Value = 10
$ws.Range("CN7").Value = 15
$ws.Range("CN8").Value = 13
$ws.Range("CN9").Value = 10
$ws.Range("CN10").Value = 16
$ws.Range("CN11").Value = 0

# Match the author's final selection recorded in the workbook view.
$ws.Range("CN2").Select()
